$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("9mm-ammo")

# Bump horizontal_recoil (E) and vertical_recoil (F) by +2 for every 9mm round (rows 11-16)
foreach ($r in 11..16) {
    $ws.Cells.Item($r, 5).Value2 = $ws.Cells.Item($r, 5).Value2 + 2
    $ws.Cells.Item($r, 6).Value2 = $ws.Cells.Item($r, 6).Value2 + 2
}

# These values were re-typed by hand, so their cell formatting reverts to the
# worksheet default (no explicit fill) instead of keeping the old style.
$ws.Range("E11:F16").Style = "Normal"

# Update conditional-formatting range to skip columns E:F (now plain cells)
$fc = $ws.Range("C11:U16").FormatConditions
$fc.Item(1).ModifyAppliesToRange($ws.Range("C11:D16,G11:U16"))

# Match the author's last on-screen selection
$ws.Range("AE23").Select()
